# Swap the "codeforiati:group-code" and "codeforiati:group-name" columns
# (columns C and D) on the active sheet: column C becomes group-name and
# column D becomes group-code for the header row and every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp

for ($r = 1; $r -le $lastRow; $r++) {
    $cCell = $ws.Cells.Item($r, 3)
    $dCell = $ws.Cells.Item($r, 4)
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2
    $cCell.Value = $dVal
    $dCell.Value = $cVal
}
